$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.01 = 7512.05 pesos`n✅ 7512.05 pesos = 2.0 = 925.02 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the numeric rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 498
$ws2.Range("O10").Value = 3741
$ws2.Range("O12").Value = 463
